$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data as captured by the diff.
# D-column (Price) cells are forced to Text format before assignment so that
# values such as "1.00", "0.375", "0.0326" are preserved as literal strings
# (matching the original inlineStr/text-typed cells) instead of being
# auto-coerced into numbers by Excel's type inference.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "91.259.25"
$ws.Range("E2").Value = "  +1.96%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.158.47"
$ws.Range("E3").Value = "  +2.52%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.53"
$ws.Range("E5").Value = "  +1.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "620.20"
$ws.Range("E6").Value = "  +0.47%  "
$ws.Range("E7").Value = "  +6.55%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.375"
$ws.Range("E8").Value = "  +3.99%  "
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.156.89"
$ws.Range("E10").Value = "  +2.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.746"
$ws.Range("E11").Value = "  +5.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.203"
$ws.Range("E12").Value = "  +2.52%  "
$ws.Range("E13").Value = "  -0.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.37"
$ws.Range("E14").Value = "  +1.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.61"
$ws.Range("E15").Value = "  +4.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.141.28"
$ws.Range("E16").Value = "  +2.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.746.15"
$ws.Range("E17").Value = "  +2.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.174.83"
$ws.Range("E18").Value = "  +3.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.76"
$ws.Range("E19").Value = "  +1.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.25"
$ws.Range("E20").Value = "  +11.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.04"
$ws.Range("E21").Value = "  +12.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "455.70"
$ws.Range("E22").Value = "  +5.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000205"
$ws.Range("E23").Value = "  -4.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.24"
$ws.Range("E24").Value = "  +6.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.03"
$ws.Range("E25").Value = "  +8.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "89.24"
$ws.Range("E26").Value = "  +3.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.08"
$ws.Range("E27").Value = "  +4.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.324.76"
$ws.Range("E28").Value = "  +2.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.133"
$ws.Range("E30").Value = "  +48.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.233"
$ws.Range("E31").Value = "  +17.54%  "
$ws.Range("E32").Value = "  +10.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.43"
$ws.Range("E33").Value = "  +4.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.172"
$ws.Range("E34").Value = "  +14.51%  "
$ws.Range("E35").Value = "  -10.63%  "
$ws.Range("E36").Value = "  +9.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.56"
$ws.Range("E37").Value = "  +3.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "513.37"
$ws.Range("E38").Value = "  +4.64%  "
$ws.Range("E39").Value = "  +5.13%  "
$ws.Range("E40").Value = "  +9.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.90"
$ws.Range("E41").Value = "  -1.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.452"
$ws.Range("E42").Value = "  +14.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.47"
$ws.Range("E43").Value = "  -3.97%  "
$ws.Range("E44").Value = "  +0.24%  "
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.718"
$ws.Range("E46").Value = "  +7.07%  "
$ws.Range("E47").Value = "  +5.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "158.40"
$ws.Range("E48").Value = "  +1.38%  "
$ws.Range("E49").Value = "  +6.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.54"
$ws.Range("E50").Value = "  +5.18%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0326"
$ws.Range("E51").Value = "  +14.57%  "
